$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 6

# Copy formatting (including the date number format / style) from the cell
# above it, then set the value - this keeps the same shared style index
# that Excel would normally reuse for matching formats.
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42588.471597222226

$ws.Cells.Item($row, 2).Value = "Random"
$ws.Cells.Item($row, 3).Value = 0
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 93
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 2
$ws.Cells.Item($row, 13).Value = 98
